$d = $word.ActiveDocument

function Replace-Text($find, $replace) {
    $d.Content.Find.Execute($find, $true, $false, $false, $false, $false, `
                             $true, 1, $false, $replace, 2)
}

function Replace-ParaKeepingEmptyRun($paraIndex, $newText, $pPrXml, $rPrXml) {
    $p = $d.Paragraphs.Item($paraIndex)
    $full = $p.Range
    $xml = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' + `
           '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' + `
           '<pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' + `
           '<w:body><w:p>' + $pPrXml + '<w:r/><w:r>' + $rPrXml + '<w:t xml:space="preserve">' + $newText + '</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
    $full.InsertXML($xml)
}

# 1. Heading1 title (no adjacent empty run, simple Find/Replace is fine) -
#    appears twice in the document (heading + bold run near the end handled separately below).
Replace-Text "Play Danger High Voltage Megapays for Free - Exciting Features and Jackpots" "Play Danger High Voltage Megapays for Free"

$listBulletPPr = '<w:pPr><w:pStyle w:val="ListBullet"/><w:spacing w:line="240" w:lineRule="auto"/><w:ind w:left="720"/></w:pPr>'

# 2-7: "What we like" / "What we don't like" bullet paragraphs (each has a leading empty run)
Replace-ParaKeepingEmptyRun 44 "Progressive jackpots" $listBulletPPr ""
Replace-ParaKeepingEmptyRun 45 "Exciting gameplay features" $listBulletPPr ""
Replace-ParaKeepingEmptyRun 46 "High volatility for big prizes" $listBulletPPr ""
Replace-ParaKeepingEmptyRun 47 "Well-made design and soundtrack" $listBulletPPr ""
Replace-ParaKeepingEmptyRun 49 "Patience required to win big" $listBulletPPr ""
Replace-ParaKeepingEmptyRun 50 "Limited bonus game options" $listBulletPPr ""

# 8: Bold "Play Danger ..." paragraph near the end (leading empty run, bold run formatting)
Replace-ParaKeepingEmptyRun 51 "Play Danger High Voltage Megapays for Free" "" "<w:rPr><w:b/></w:rPr>"

# 9: Italic summary paragraph at the very end (leading empty run, italic run formatting)
Replace-ParaKeepingEmptyRun 52 "Read our review of Danger High Voltage Megapays and play for free on your mobile device." "" "<w:rPr><w:i/></w:rPr>"
